# Update profit files after running on 2025-11-02
# Appends the new daily allocation row (11/02/2025) to the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 62

# Force the date column to be stored as literal text (matching the rest of
# column A, which holds "MM/DD/YYYY" strings rather than real date values),
# then drop the temporary number format so the cell keeps the sheet's
# default (unstyled) look.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "11/02/2025"
$dateCell.ClearFormats()

$ws.Cells.Item($row, 2).Value = 0.1979663639690586
$ws.Cells.Item($row, 3).Value = 0.8020336360309414
